# Fill in previously-empty benchmark cells (inlineStr cells with no <is>) on the
# "BENCHMARK" worksheet, rows 3-14, columns D/F/G/I, matching the updated
# benchmark figures from the 2025-09-14 13:24:13 UTC run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - HESAPTAN EFT - Şube
$ws.Range("D3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("G3").Value = "30,46 TRY - 60,94 TRY - 609,43 TRY"
$ws.Range("I3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 4 - HESAPTAN EFT - ATM
$ws.Range("D4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("G4").Value = "21,27 TRY - 42,55 TRY - 304,71 TRY"
$ws.Range("I4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 5 - HESAPTAN EFT - Mobil
$ws.Range("D5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("G5").Value = "6,09 TRY - 12,19 TRY - 152,35 TRY"
$ws.Range("I5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 6 - DÜZENLİ EFT  (F6 stays empty)
$ws.Range("D6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("G6").Value = "4.300,01 TL - 76,17 TL"
$ws.Range("I6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

# Row 8 - HESAPTAN HAVALE - Şube
$ws.Range("D8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("G8").Value = "15,23 TRY - 30,47 TRY - 304,71 TRY"
$ws.Range("I8").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 9 - HESAPTAN HAVALE - ATM
$ws.Range("D9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("G9").Value = "10,63 TRY - 21,27 TRY - 152,35 TRY"
$ws.Range("I9").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 10 - HESAPTAN HAVALE - Mobil
$ws.Range("D10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("G10").Value = "3,04 TRY - 6,09 TRY - 76,17 TRY"
$ws.Range("I10").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 11 - DÜZENLİ HAVALE  (F11 stays empty)
$ws.Range("D11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("G11").Value = "3,04 TRY - 6,09 TRY - 76,17 TRY"
$ws.Range("I11").Value = "3,04 TL - 6,09 TL - 76,17 TL"

# Row 12 - GİDEN SWIFT  (E12, F12, H12, I12, J12 stay empty)
$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"
$ws.Range("G12").Value = "Şube (Kasadan): %0,5; Şube (Hesaptan): %0,75; İnternet: 15 USD"

# Row 13 - GELEN SWIFT  (G13 stays empty)
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("F13").Value = "Hesaba: Asgari 300 TL | Azami 3.080 TL"
$ws.Range("I13").Value = "Hesaba: Asgari 1 TL | Azami 6,09 TL"

# Row 14 - GİDEN SWIFT - Mobil
$ws.Range("D14").Value = "2.300 TL - 9.500 TL"
$ws.Range("F14").Value = "1.952,38 TL - 9.523,81 TL"
$ws.Range("G14").Value = "4.300 TL - 6,09 TL"
